$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.819.84"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "1.891.64"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7742"
$ws.Range("E5").Value = "  -4.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.07"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3128"
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.35"
$ws.Range("E9").Value = "  -7.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07227"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08080"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7654"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "1.913.59"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.35"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.150"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "29.845.51"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.94"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.53"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007764"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "2.153.76"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.115"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1564"
$ws.Range("E25").Value = "  -3.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.396"
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.40"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.044"
$ws.Range("E29").Value = "  -4.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.439"
$ws.Range("E30").Value = "  +4.58%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.475"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.101"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05498"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.260"
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7461"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.632"
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01917"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.787"
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").Value = "1.140.37"
$ws.Range("E41").Value = "  +10.71%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.886"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8496"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.96"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.885"
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.893"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.044"
$ws.Range("E50").Value = "  +11.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.455"
$ws.Range("E51").Value = "  -2.89%  "
